# Adds two new columns ("BTTS?" and "BTTS p?") to the score-prediction
# sheet, fills in their values, extends the summary SUM / percentage rows,
# and mirrors the header formatting + a highlight fill on the new total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column data (rows 2-49) ------------------------------------------
# Column M = "BTTS?" (both teams scored - actual)
# Column L = "BTTS p?" (both teams scored - predicted)
$bttsPredicted = @(1,1,1,0,0,1,1,0,1,1,0,1,1,1,1,1,1,1,1,0,1,1,1,0,1,1,1,0,1,1,1,1,1,1,0,1,1,1,1,1,1,1,0,1,1,0,1,1)
$bttsActual    = @(0,1,1,0,1,0,0,1,1,1,1,0,0,1,0,0,0,1,0,0,1,0,1,0,1,1,1,0,1,1,0,0,1,0,0,1,1,1,0,0,0,1,0,0,0,0,1,1)

# --- Headers ----------------------------------------------------------------
# Write M1 before L1 so the new shared-string entries land in the same order
# as the source workbook (index 58 "BTTS?", index 59 "BTTS p?").
$ws.Cells.Item(1, 13).Value = "BTTS?"
$ws.Cells.Item(1, 12).Value = "BTTS p?"

# Copy the existing header look (bold font, right border, centered/top) from
# K1 onto the two new header cells.
$ws.Range("K1").Copy()
$ws.Range("L1:M1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row values --------------------------------------------------------------
for ($i = 0; $i -lt 48; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 12).Value = $bttsPredicted[$i]
    $ws.Cells.Item($row, 13).Value = $bttsActual[$i]
}

# --- Totals row (50): extend the SUM() pattern already used in K50 ----------
$ws.Range("L50").Formula = "=SUM(L2:L49)"
$ws.Range("M50").Formula = "=SUM(M2:M49)"

# Highlight the predicted-BTTS total with the new orange fill.
$ws.Range("L50").Interior.Color = 49407

# --- Percentage row (51): extend the 100*(x/48) pattern already used -------
$ws.Range("L51").Formula = "=100*(L50/48)"
$ws.Range("M51").Formula = "=100*(M50/48)"

# --- Selection / active cell, matching where editing finished --------------
$ws.Range("A6").Select()
$ws.Range("L50").Select()
